$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 1.440985666666667
$ws.Range("H2").Value = 4.322957000000001
$ws.Range("I2").Value = 0.1098365531732288
$ws.Range("J2").Value = 0.1230162332390494
$ws.Range("M2").Value = 28.85518433333334
$ws.Range("N2").Value = 86.56555300000001
$ws.Range("O2").Value = 0.1999651185353207
$ws.Range("P2").Value = 0.2044513327926365
$ws.Range("Q2").Value = 41.57990703335791
$ws.Range("R2").Value = 374.2191633002211
$ws.Range("S2").Value = 0.02196347937479575
$ws.Range("T2").Value = 0.02515083284085348
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 1.440985666666667
$ws.Range("H3").Value = 4.322957000000001
$ws.Range("I3").Value = 0.1098365531732288
$ws.Range("J3").Value = 0.1230162332390494
$ws.Range("O3").Value = 0.3546352265743414
$ws.Range("P3").Value = 0.3625914622481308
$ws.Range("Q3").Value = 73.74135979175958
$ws.Range("R3").Value = 663.6722381258361
$ws.Range("S3").Value = 0.0389519109207327
$ws.Range("T3").Value = 0.04460463589040402
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 1.440985666666667
$ws.Range("H4").Value = 4.322957000000001
$ws.Range("I4").Value = 0.1098365531732288
$ws.Range("J4").Value = 0.1230162332390494
$ws.Range("M4").Value = 29.393479
$ws.Range("N4").Value = 88.180437
$ws.Range("O4").Value = 0.2036954761578358
$ws.Range("P4").Value = 0.2082653809291453
$ws.Range("Q4").Value = 42.35558193246767
$ws.Range("R4").Value = 381.200237392209
$ws.Range("S4").Value = 0.0223732089981563
$ws.Range("T4").Value = 0.0256200226759992
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 1.440985666666667
$ws.Range("H5").Value = 4.322957000000001
$ws.Range("I5").Value = 0.1098365531732288
$ws.Range("J5").Value = 0.1230162332390494
$ws.Range("M5").Value = 9.499066500000001
$ws.Range("N5").Value = 18.998133
$ws.Range("O5").Value = 0.0658280999596015
$ws.Range("P5").Value = 0.04486996822421697
$ws.Range("Q5").Value = 13.6880186732135
$ws.Range("R5").Value = 82.12811203928102
$ws.Range("S5").Value = 0.00723033160150539
$ws.Range("T5").Value = 0.005519734476499008
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 1.440985666666667
$ws.Range("H6").Value = 4.322957000000001
$ws.Range("I6").Value = 0.1098365531732288
$ws.Range("J6").Value = 0.1230162332390494
$ws.Range("M6").Value = 25.37910966666666
$ws.Range("N6").Value = 76.13732899999999
$ws.Range("O6").Value = 0.1758760787729007
$ws.Range("P6").Value = 0.1798218558058706
$ws.Range("Q6").Value = 36.57093326242811
$ws.Range("R6").Value = 329.138399361853
$ws.Range("S6").Value = 0.01931762227803868
$ws.Range("T6").Value = 0.02212100735529368
$ws.Range("I7").Value = 0.5687502547919595
$ws.Range("J7").Value = 0.6369966279614609
$ws.Range("M7").Value = 28.85518433333334
$ws.Range("N7").Value = 86.56555300000001
$ws.Range("O7").Value = 0.1999651185353207
$ws.Range("P7").Value = 0.2044513327926365
$ws.Range("Q7").Value = 215.3070361025524
$ws.Range("R7").Value = 1937.763324922972
$ws.Range("S7").Value = 0.113730212116468
$ws.Range("T7").Value = 0.1302348095711359
$ws.Range("I8").Value = 0.5687502547919595
$ws.Range("J8").Value = 0.6369966279614609
$ws.Range("O8").Value = 0.3546352265743414
$ws.Range("P8").Value = 0.3625914622481308
$ws.Range("S8").Value = 0.201698875472361
$ws.Range("T8").Value = 0.2309695387796747
$ws.Range("I9").Value = 0.5687502547919595
$ws.Range("J9").Value = 0.6369966279614609
$ws.Range("M9").Value = 29.393479
$ws.Range("N9").Value = 88.180437
$ws.Range("O9").Value = 0.2036954761578358
$ws.Range("P9").Value = 0.2082653809291453
$ws.Range("Q9").Value = 219.3235978368653
$ws.Range("R9").Value = 1973.912380531788
$ws.Range("S9").Value = 0.1158518539647387
$ws.Range("T9").Value = 0.1326643453729747
$ws.Range("I10").Value = 0.5687502547919595
$ws.Range("J10").Value = 0.6369966279614609
$ws.Range("M10").Value = 9.499066500000001
$ws.Range("N10").Value = 18.998133
$ws.Range("O10").Value = 0.0658280999596015
$ws.Range("P10").Value = 0.04486996822421697
$ws.Range("Q10").Value = 70.878627224482
$ws.Range("R10").Value = 425.271763346892
$ws.Range("S10").Value = 0.03743974862449393
$ws.Range("T10").Value = 0.02858201845556411
$ws.Range("I11").Value = 0.5687502547919595
$ws.Range("J11").Value = 0.6369966279614609
$ws.Range("M11").Value = 25.37910966666666
$ws.Range("N11").Value = 76.13732899999999
$ws.Range("O11").Value = 0.1758760787729007
$ws.Range("P11").Value = 0.1798218558058706
$ws.Range("Q11").Value = 189.3698136919995
$ws.Range("R11").Value = 1704.328323227996
$ws.Range("S11").Value = 0.100029564613898
$ws.Range("T11").Value = 0.1145459157821117
$ws.Range("G12").Value = 4.2167365
$ws.Range("H12").Value = 8.433472999999999
$ws.Range("I12").Value = 0.3214131920348118
$ws.Range("J12").Value = 0.2399871387994896
$ws.Range("M12").Value = 28.85518433333334
$ws.Range("N12").Value = 86.56555300000001
$ws.Range("O12").Value = 0.1999651185353207
$ws.Range("P12").Value = 0.2044513327926365
$ws.Range("Q12").Value = 121.6747089925948
$ws.Range("R12").Value = 730.0482539555691
$ws.Range("S12").Value = 0.06427142704405693
$ws.Range("T12").Value = 0.04906569038064709
$ws.Range("G13").Value = 4.2167365
$ws.Range("H13").Value = 8.433472999999999
$ws.Range("I13").Value = 0.3214131920348118
$ws.Range("J13").Value = 0.2399871387994896
$ws.Range("O13").Value = 0.3546352265743414
$ws.Range("P13").Value = 0.3625914622481308
$ws.Range("Q13").Value = 215.7883250239673
$ws.Range("R13").Value = 1294.729950143804
$ws.Range("S13").Value = 0.1139844401812478
$ws.Range("T13").Value = 0.08701728757805206
$ws.Range("G14").Value = 4.2167365
$ws.Range("H14").Value = 8.433472999999999
$ws.Range("I14").Value = 0.3214131920348118
$ws.Range("J14").Value = 0.2399871387994896
$ws.Range("M14").Value = 29.393479
$ws.Range("N14").Value = 88.180437
$ws.Range("O14").Value = 0.2036954761578358
$ws.Range("P14").Value = 0.2082653809291453
$ws.Range("Q14").Value = 123.9445557612835
$ws.Range("R14").Value = 743.6673345677009
$ws.Range("S14").Value = 0.06547041319494092
$ws.Range("T14").Value = 0.04998101288017136
$ws.Range("G15").Value = 4.2167365
$ws.Range("H15").Value = 8.433472999999999
$ws.Range("I15").Value = 0.3214131920348118
$ws.Range("J15").Value = 0.2399871387994896
$ws.Range("M15").Value = 9.499066500000001
$ws.Range("N15").Value = 18.998133
$ws.Range("O15").Value = 0.0658280999596015
$ws.Range("P15").Value = 0.04486996822421697
$ws.Range("Q15").Value = 40.05506042647725
$ws.Range("R15").Value = 160.220241705909
$ws.Range("S15").Value = 0.02115801973360218
$ws.Range("T15").Value = 0.01076821529215385
$ws.Range("G16").Value = 4.2167365
$ws.Range("H16").Value = 8.433472999999999
$ws.Range("I16").Value = 0.3214131920348118
$ws.Range("J16").Value = 0.2399871387994896
$ws.Range("M16").Value = 25.37910966666666
$ws.Range("N16").Value = 76.13732899999999
$ws.Range("O16").Value = 0.1758760787729007
$ws.Range("P16").Value = 0.1798218558058706
$ws.Range("Q16").Value = 107.0170180689362
$ws.Range("R16").Value = 642.1021084136169
$ws.Range("S16").Value = 0.056528891880964
$ws.Range("T16").Value = 0.04315493266846528
